# Generate Report for Handback
#
# "befd76d3-acd3-492d-aa12-05c8a914ab23" has been handed back (translation
# complete, in sync with en-US) while "d3e7fab8-0350-4ca7-86dd-e36c0d26afb4"
# is still in translation. The status report rows for these two files swap
# places (alphabetical ordering by file name) and the befd76d3 row picks up
# its handback info (target file / handback file / handback datetime).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "Overview" sheet - rows 5 and 6 swap (file name + status columns)
# ---------------------------------------------------------------------
$ovw = $wb.Worksheets.Item("Overview")

$ovw.Range("A5").Value = "befd76d3-acd3-492d-aa12-05c8a914ab23.md"
$ovw.Range("B5").Value = "Handed back: in sync with en-US"
$ovw.Range("C5").Value = "Handed back: in sync with en-US"

$ovw.Range("A6").Value = "d3e7fab8-0350-4ca7-86dd-e36c0d26afb4.md"
$ovw.Range("B6").Value = "In Translation"
$ovw.Range("C6").Value = "In Translation"

# ---------------------------------------------------------------------
# "zh-cn" sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

# Row 5 becomes the befd76d3 entry (now handed back)
$zh.Range("A5").Value = "befd76d3-acd3-492d-aa12-05c8a914ab23.md"
$zh.Range("B5").Value = "Handed back: in sync with en-US"
$zh.Range("C5").Value = "befd76d3-acd3-492d-aa12-05c8a914ab23.185d1e5cbeb72050f7f5fa46bc372a721eeb1ae2.zh-cn.xlf"
$zh.Range("D5").Value = "2016-03-11 02:24:50"
$zh.Range("E5").Value = "befd76d3-acd3-492d-aa12-05c8a914ab23.md"
$zh.Range("F5").Value = "befd76d3-acd3-492d-aa12-05c8a914ab23.185d1e5cbeb72050f7f5fa46bc372a721eeb1ae2.zh-cn.xlf"
$zh.Range("G5").Value = "2016-03-11 02:25:44"
# give the newly populated hyperlink-like cells the same look as their neighbours
$zh.Range("E5").Font.Name = "Calibri"
$zh.Range("E5").Font.Size = 11
$zh.Range("E5").Font.Underline = 2
$zh.Range("E5").Font.Color = 15570276
$zh.Range("F5").Font.Name = "Calibri"
$zh.Range("F5").Font.Size = 11
$zh.Range("F5").Font.Underline = 2
$zh.Range("F5").Font.Color = 15570276

# Row 6 becomes the d3e7fab8 entry (still in translation, no handback yet)
$zh.Range("A6").Value = "d3e7fab8-0350-4ca7-86dd-e36c0d26afb4.md"
$zh.Range("B6").Value = "In Translation"
$zh.Range("C6").Value = "d3e7fab8-0350-4ca7-86dd-e36c0d26afb4.a9567d8361ef552a0252e4f39417c927a83e4a86.zh-cn.xlf"
$zh.Range("D6").Value = "2016-03-11 02:15:37"
$zh.Range("E6").Value = ""
$zh.Range("F6").Value = ""
$zh.Range("G6").Value = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# "de-de" sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

# Row 5 becomes the befd76d3 entry (now handed back)
$de.Range("A5").Value = "befd76d3-acd3-492d-aa12-05c8a914ab23.md"
$de.Range("B5").Value = "Handed back: in sync with en-US"
$de.Range("C5").Value = "befd76d3-acd3-492d-aa12-05c8a914ab23.185d1e5cbeb72050f7f5fa46bc372a721eeb1ae2.de-de.xlf"
$de.Range("D5").Value = "2016-03-11 02:24:58"
$de.Range("E5").Value = "befd76d3-acd3-492d-aa12-05c8a914ab23.md"
$de.Range("F5").Value = "befd76d3-acd3-492d-aa12-05c8a914ab23.185d1e5cbeb72050f7f5fa46bc372a721eeb1ae2.de-de.xlf"
$de.Range("G5").Value = "2016-03-11 02:26:07"
$de.Range("E5").Font.Name = "Calibri"
$de.Range("E5").Font.Size = 11
$de.Range("E5").Font.Underline = 2
$de.Range("E5").Font.Color = 15570276
$de.Range("F5").Font.Name = "Calibri"
$de.Range("F5").Font.Size = 11
$de.Range("F5").Font.Underline = 2
$de.Range("F5").Font.Color = 15570276

# Row 6 becomes the d3e7fab8 entry (still in translation, no handback yet)
$de.Range("A6").Value = "d3e7fab8-0350-4ca7-86dd-e36c0d26afb4.md"
$de.Range("B6").Value = "In Translation"
$de.Range("C6").Value = "d3e7fab8-0350-4ca7-86dd-e36c0d26afb4.a9567d8361ef552a0252e4f39417c927a83e4a86.de-de.xlf"
$de.Range("D6").Value = "2016-03-11 02:17:03"
$de.Range("E6").Value = ""
$de.Range("F6").Value = ""
$de.Range("G6").Value = "0001-01-01 00:00:00"
